$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers; force them to stay text
$textCellAddrs = @("D8", "D11", "D16", "D18", "D20", "D22", "D25", "D26", "D28", "D34", "D38", "D39", "D43", "D46", "D48", "D51")
foreach ($addr in $textCellAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '27.961.39'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.638.51'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '23.31'
$ws.Range('E8').Value = '  -1.51%  '
$ws.Range('E9').Value = '  -2.44%  '
$ws.Range('E10').Value = '  -0.01%  '
$ws.Range('D11').Value = '0.0883'
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').Value = '1.870.77'
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').Value = '1.639.78'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('E15').Value = '  +0.84%  '
$ws.Range('D16').Value = '65.34'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('D17').Value = '27.963.03'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '231.39'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').Value = '7.55'
$ws.Range('E20').Value = '  -1.99%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '10.43'
$ws.Range('E22').Value = '  -2.48%  '
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('E24').Value = '  -3.88%  '
$ws.Range('D25').Value = '153.56'
$ws.Range('E25').Value = '  +1.25%  '
$ws.Range('D26').Value = '6.98'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').Value = '15.64'
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  -0.40%  '
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').Value = '1.405.72'
$ws.Range('E33').Value = '  -3.50%  '
$ws.Range('D34').Value = '3.08'
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('E35').Value = '  +1.31%  '
$ws.Range('E36').Value = '  +1.58%  '
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = '0.563'
$ws.Range('E38').Value = '  -0.22%  '
$ws.Range('D39').Value = '0.926'
$ws.Range('E39').Value = '  +0.33%  '
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '67.02'
$ws.Range('E43').Value = '  -3.62%  '
$ws.Range('E44').Value = '  +2.58%  '
$ws.Range('E45').Value = '  +1.64%  '
$ws.Range('D46').Value = '2.21'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').Value = '1.780.33'
$ws.Range('E47').Value = '  -0.54%  '
$ws.Range('D48').Value = '88.10'
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('E49').Value = '  -0.72%  '
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').Value = '7.56'
$ws.Range('E51').Value = '  -2.00%  '

# Restore original (default/General) style now that text values are locked in
foreach ($addr in $textCellAddrs) {
    $ws.Range($addr).Style = "Normal"
}

